$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45985
$ws.Range("B2").Value = 27.6
$ws.Range("C2").Value = 10.8
$ws.Range("D2").Value = 3.91
$ws.Range("E2").Value = 3.52
$ws.Range("F2").Value = 3.52
$ws.Range("G2").Value = 3.52
$ws.Range("H2").Value = 15.94
$ws.Range("I2").Value = 47.08
$ws.Range("J2").Value = 61.28
$ws.Range("K2").Value = 47.88
$ws.Range("L2").Value = 39.2
$ws.Range("M2").Value = 31.49
$ws.Range("N2").Value = 28.13
$ws.Range("O2").Value = 34.17
$ws.Range("P2").Value = 38.67
$ws.Range("Q2").Value = 45.26
$ws.Range("R2").Value = 55.79
$ws.Range("S2").Value = 61.68
$ws.Range("T2").Value = 73.06
$ws.Range("U2").Value = 74.82
$ws.Range("V2").Value = 82.07
$ws.Range("W2").Value = 77.74
$ws.Range("X2").Value = 71.27
$ws.Range("Y2").Value = 64.2
$ws.Range("Z2").Value = 41.78
$ws.Range("AB2").Value = 73.82
$ws.Range("AD2").Value = 79.9
$ws.Range("AF2").Value = 73.94
$ws.Range("AG2").Value = "0h-14h"
